# Delete slide 18 ("Ejemplo de Diagrama de Secuencia 1"), which merges
# the duplicate "Diagrama de Secuencia" slides down to a single one
# ("Ejemplo de Diagrama de Secuencia 2") and shifts all following
# slides up by one position.
$p = $ppt.ActivePresentation
$p.Slides.Item(18).Delete()
